# Update cryptocurrency price/volume figures per the latest crawl snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.766.44'
$ws.Range('E2').Value = '  +1.78%  '
$ws.Range('D3').Value = '2.093.05'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.39'
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.21'
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0846'
$ws.Range('E10').Value = '  +0.61%  '
$ws.Range('D11').Value = '3.011.05'
$ws.Range('E11').Value = '  +25.56%  '
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.43'
$ws.Range('E13').Value = '  +5.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.12'
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.809'
$ws.Range('E15').Value = '  +4.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.47'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('D17').Value = '2.094.70'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '38.745.34'
$ws.Range('E18').Value = '  +3.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.88'
$ws.Range('E19').Value = '  +2.63%  '
$ws.Range('E20').Value = '  +0.67%  '
$ws.Range('D21').Value = '0.0₃0842'
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.70'
$ws.Range('E22').Value = '  +1.57%  '
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.38'
$ws.Range('E24').Value = '  -1.93%  '
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.30'
$ws.Range('E26').Value = '  +0.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.56'
$ws.Range('E27').Value = '  +1.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.140'
$ws.Range('E28').Value = '  +5.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.30'
$ws.Range('E30').Value = '  +1.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.49'
$ws.Range('E31').Value = '  +5.31%  '
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.52'
$ws.Range('E33').Value = '  +1.79%  '
$ws.Range('E34').Value = '  +0.80%  '
$ws.Range('E35').Value = '  +0.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.57'
$ws.Range('E36').Value = '  +1.52%  '
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.61'
$ws.Range('E38').Value = '  +1.90%  '
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.02'
$ws.Range('E40').Value = '  +0.13%  '
$ws.Range('E41').Value = '  +4.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '101.15'
$ws.Range('E42').Value = '  +1.19%  '
$ws.Range('D43').Value = '1.533.04'
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.80'
$ws.Range('E44').Value = '  -0.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0915'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.14'
$ws.Range('E46').Value = '  +1.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.66'
$ws.Range('E47').Value = '  +5.81%  '
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.04'
$ws.Range('E49').Value = '  +1.19%  '
$ws.Range('E50').Value = '  -0.93%  '
$ws.Range('D51').Value = '2.287.67'
$ws.Range('E51').Value = '  +0.07%  '
